# The deck's single Design (slide master) currently uses the "Integral" /
# "Red Violet" colour theme (stored in ppt/theme/theme1.xml). The target
# revision swaps that in for the stock "Office Theme" colour palette (the
# palette that, in this deck, otherwise only lived in the notes-master
# theme part). Re-apply every theme colour on the active design's theme so
# the slide master (and every slide/layout that inherits from it) picks up
# the Office Theme palette.
#
# PowerPoint's ThemeColorScheme colours are standard OLE RGB integers
# (0xBBGGRR = B*65536 + G*256 + R), so each target hex colour below is
# converted accordingly:
#   dk1      000000 -> 0
#   lt1      FFFFFF -> 16777215
#   dk2      44546A -> 6968388
#   lt2      E7E6E6 -> 15132391
#   accent1  5B9BD5 -> 13998939
#   accent2  ED7D31 -> 3243501
#   accent3  A5A5A5 -> 10855845
#   accent4  FFC000 -> 49407
#   accent5  4472C4 -> 12874308
#   accent6  70AD47 -> 4697456
#   hlink    0563C1 -> 12673797
#   folHlink 954F72 -> 7491477

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
